$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.621.33'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '1.840.94'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'259.29"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'0.5289"
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("D8").Value = "'0.3136"
$ws.Range("E8").Value = '  -3.71%  '
$ws.Range("D9").Value = "'0.06799"
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("D10").Value = "'18.69"
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("D11").Value = "'0.7790"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = "'0.07760"
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '1.830.34'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = "'87.86"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = "'5.008"
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("D16").Value = "'0.9994"
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = "'13.82"
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = "'1.0000"
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = "'0.000007919"
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '26.640.52'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").Value = '2.077.12'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = "'4.599"
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = "'5.969"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = "'9.315"
$ws.Range("E24").Value = '  -2.23%  '
$ws.Range("D25").Value = "'142.91"
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("D26").Value = "'2.207"
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").Value = "'1.684"
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("D28").Value = "'16.97"
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = "'110.89"
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").Value = "'4.181"
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = "'0.08730"
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").Value = "'4.066"
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("D36").Value = "'2.862"
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("D37").Value = "'3.089"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = "'2.250"
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("D39").Value = "'0.01726"
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").Value = "'0.4792"
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").Value = "'0.8953"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").Value = "'109.84"
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").Value = "'5.916"
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = "'7.664"
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").Value = "'0.4157"
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("D47").Value = "'8.967"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = "'0.1236"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").Value = "'0.05820"
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").Value = "'34.71"
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").Value = "'0.8934"
$ws.Range("E51").Value = '  +0.86%  '
